$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1), matching the formatting of the
# existing header row (copy style from G1 so the same cellXf is reused).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New "Save" column data values (H2:H5)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
